$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Learnings of 16th Nov" entries: dates for Nov 12 - Nov 16, 2020
$dates = @(44147, 44148, 44149, 44150, 44151)
$startRow = 14

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
}

# Copy the existing date-formatted style (numFmtId 14) from A13 onto the
# new date cells so we reuse the same cellXfs entry instead of creating one.
$ws.Range("A13").Copy()
$ws.Range("A14:A15").PasteSpecial(-4122)
$ws.Range("A17:A18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 16's date cell ended up with a different ("d-mmm") number format.
$ws.Cells.Item(16, 1).NumberFormat = "d-mmm"

# Topic learned on the 16th
$ws.Cells.Item(18, 2).Value = "JSON"

$ws.Range("B18").Select()
